$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.08840121110588733
$ws.Range("J2").Value = 0.08840121110588735
$ws.Range("M2").Value = 26.34406266666667
$ws.Range("N2").Value = 79.032188
$ws.Range("O2").Value = 0.3168347904490542
$ws.Range("P2").Value = 0.3168347904490542
$ws.Range("Q2").Value = 9.378670717772001
$ws.Range("R2").Value = 84.40803645994801
$ws.Range("S2").Value = 0.02800857919617642
$ws.Range("T2").Value = 0.02800857919617642

# Row 3
$ws.Range("I3").Value = 0.08840121110588733
$ws.Range("J3").Value = 0.08840121110588735
$ws.Range("O3").Value = 0.0001452797431229321
$ws.Range("P3").Value = 0.0001452797431229321
$ws.Range("Q3").Value = 0.004300445891
$ws.Range("R3").Value = 0.03870401301900001
$ws.Range("S3").Value = [double]"1.28429052412194E-05"
$ws.Range("T3").Value = [double]"1.28429052412194E-05"

# Row 4
$ws.Range("I4").Value = 0.08840121110588733
$ws.Range("J4").Value = 0.08840121110588735
$ws.Range("O4").Value = 0.6830199298078229
$ws.Range("P4").Value = 0.6830199298078229
$ws.Range("S4").Value = 0.0603797890044697
$ws.Range("T4").Value = 0.06037978900446971

# Row 5
$ws.Range("M5").Value = 26.34406266666667
$ws.Range("N5").Value = 79.032188
$ws.Range("O5").Value = 0.3168347904490542
$ws.Range("P5").Value = 0.3168347904490542
$ws.Range("Q5").Value = 96.71343594509023
$ws.Range("R5").Value = 870.420923505812
$ws.Range("S5").Value = 0.2888262112528778
$ws.Range("T5").Value = 0.2888262112528778

# Row 6
$ws.Range("O6").Value = 0.0001452797431229321
$ws.Range("P6").Value = 0.0001452797431229321
$ws.Range("S6").Value = 0.0001324368378817127
$ws.Range("T6").Value = 0.0001324368378817127

# Row 7
$ws.Range("O7").Value = 0.6830199298078229
$ws.Range("P7").Value = 0.6830199298078229
$ws.Range("S7").Value = 0.6226401408033532
$ws.Range("T7").Value = 0.6226401408033533
